$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# Correct the Programme / Cost Centre values on row 2
$ws.Range("A2").Value = 109076
$ws.Range("B2").Value = 11272001

# Size column B to fit its (now numeric) contents, as Excel's "best fit"
# would after the data correction above.
$ws.Columns.Item(2).ColumnWidth = 8.3

# The active selection collapses from the whole row to just A4.
$ws.Range("A4").Select() | Out-Null
